# Extend the "stg_card" template (Sheet1) with a new "desc" (string) column.
# This mirrors the existing id/name/.../model_res columns: row1 = field name,
# row2 = field type, row3 = Chinese description, row4 = sample data (left
# blank for this new column, same as in the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G1").Value = "desc"
$ws.Range("G2").Value = "string"
$ws.Range("G3").Value = "描述"

# Move/leave the active selection on the newly added description cell,
# matching the authored selection change (G7 -> G3).
$ws.Range("G3").Select()
